# Updated time stamp illustration
# Applies the geometry / adjustment / font-size changes captured in the
# target OOXML diff for slide 5 of the deck.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# ---------------------------------------------------------------
# 1. The three arc shapes (116 / 117 / 118)
# ---------------------------------------------------------------
$sh116 = Get-ShapeById $s.Shapes 116
$sh116.Left   = 404.7044094488189
$sh116.Top    = 226.89574803149605
$sh116.Width  = 189.12496062992125
$sh116.Height = 234.92700787401574
$sh116.Adjustments.Item(1) = 153.47787

$sh117 = Get-ShapeById $s.Shapes 117
$sh117.Left   = 407.31307086614174
$sh117.Top    = 239.0708661417323
$sh117.Width  = 144.23874015748032
$sh117.Height = 166.58992125984253
$sh117.Adjustments.Item(1) = 163.88594

$sh118 = Get-ShapeById $s.Shapes 118
$sh118.Left   = 283.6827559055118
$sh118.Top    = 105.12141732283465
$sh118.Width  = 200.8983464566929
$sh118.Height = 204.09629921259844
$sh118.Adjustments.Item(2) = 55.25388

# ---------------------------------------------------------------
# 2. The straight connectors (120 / 121 / 122 / 125 / 126 / 151)
# ---------------------------------------------------------------
$sh120 = Get-ShapeById $s.Shapes 120
$sh120.Width  = 56.09456692913386
$sh120.Height = 67.25818897637795

$sh121 = Get-ShapeById $s.Shapes 121
$sh121.Left   = 385.18047244094487
$sh121.Top    = 240.99330708661418
$sh121.Width  = 62.62283464566929
$sh121.Height = 68.90204724409449

$sh122 = Get-ShapeById $s.Shapes 122
$sh122.Left   = 447.80330708661415
$sh122.Top    = 240.88
$sh122.Width  = 45.91771653543307
$sh122.Height = 70.17125984251969

$sh125 = Get-ShapeById $s.Shapes 125
$sh125.Left   = 493.72094488188975
$sh125.Top    = 310.4607874015748
$sh125.Width  = 40.81732283464567
$sh125.Height = 77.64818897637795

$sh126 = Get-ShapeById $s.Shapes 126
$sh126.Left   = 439.54527559055117
$sh126.Top    = 384.83314960629923
$sh126.Width  = 43.4144094488189
$sh126.Height = 72.68244094488189

$sh151 = Get-ShapeById $s.Shapes 151
$sh151.Left   = 440.1748031496063
$sh151.Top    = 310.87535433070866
$sh151.Width  = 53.54606299212598
$sh151.Height = 75.22685039370079

# ---------------------------------------------------------------
# 3. The three small marker groups (272 / 302 / 322) only move
#    horizontally.
# ---------------------------------------------------------------
$sh272 = Get-ShapeById $s.Shapes 272
$sh272.Left = 425.8068503937008

$sh302 = Get-ShapeById $s.Shapes 302
$sh302.Left = 418.36535433070867

$sh322 = Get-ShapeById $s.Shapes 322
$sh322.Left = 511.201968503937

# ---------------------------------------------------------------
# 4. Big "Group 18" (id 19) — blue numbers — resize group envelope
#    and reposition / resize + enlarge font of its children.
# ---------------------------------------------------------------
$grp19 = Get-ShapeById $s.Shapes 19
$grp19.Width  = 228.9139
$grp19.Height = 389.1327

function Set-Num($grp, $id, $left, $top, $width, $height) {
    $c = Get-ShapeById $grp.GroupItems $id
    if ($left -ne $null)   { $c.Left = $left }
    if ($top -ne $null)    { $c.Top = $top }
    $c.Width  = $width
    $c.Height = $height
    $c.TextFrame.TextRange.Font.Size = 20
}

Set-Num $grp19 247 436.0028346456693 463.5635433070866 44.59724409448819 31.50472440944882
Set-Num $grp19 258 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp19 268 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp19 278 405.957874015748 242.6775590551181 34.23055118110236 31.50472440944882
Set-Num $grp19 288 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp19 298 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp19 308 380.56110236220474 386.84771653543305 52.31078740157481 31.50472440944882
Set-Num $grp19 318 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp19 328 487.60165354330707 386.80889763779527 38.48188976377953 31.50472440944882
Set-Num $grp19 338 $null $null 34.23055118110236 31.50472440944882

# ---------------------------------------------------------------
# 5. Big "Group 17" (id 18) — green numbers — same treatment.
# ---------------------------------------------------------------
$grp18 = Get-ShapeById $s.Shapes 18
$grp18.Width  = 232.4153543307087
$grp18.Height = 389.1327

Set-Num $grp18 250 493.7930708661417 463.5635433070866 43.04496062992126 31.50472440944882
Set-Num $grp18 256 400.20543307086615 105.93566929133858 41.49787401574803 31.50472440944882
Set-Num $grp18 266 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp18 276 453.6315748031496 242.6775590551181 48.898110236220475 31.50472440944882
Set-Num $grp18 286 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp18 296 $null $null 34.23055118110236 31.50472440944882
Set-Num $grp18 306 446.4400787401575 386.84771653543305 41.886771653543306 31.50472440944882
Set-Num $grp18 316 399.70543307086615 176.26858267716534 47.53503937007874 31.50472440944882
Set-Num $grp18 326 539.2766929133858 386.80889763779527 38.48188976377953 31.50472440944882
Set-Num $grp18 336 $null $null 42.94062992125984 31.50472440944882
